$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (pushes existing rows 14+ down to 15+),
# inheriting formatting from the row above (row 13: "Docentes responsáveis" entry).
$ws.Rows.Item(14).Insert()

# Populate the new row with the second instructor, same layout as row 13
# (label only in column A on the "header" row above; data duplicated in B/C).
$ws.Range("B14").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C14").Value = "7290967 - Emerson Gonçalves de Melo"
